# Natmi following Dr Hou advice
# Updates ligand/receptor-expressing-cell counts (E,K: 1 -> 3) and
# the dependent expression/specificity/weight statistics (G,H,I,J,M,N,O,P,Q,R,S,T)
# for every data row (2-17) in the NATMI LR-pairs sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number => @{ column letter = new value }
$rowUpdates = @{
    2 = @{ "E"=3; "G"=31.29437; "H"=93.88310999999999; "I"=0.2388439922596655; "J"=0.2388439922596655; "K"=3; "M"=3.778439; "N"=11.335317; "O"=0.4252971528324392; "P"=0.4252971528324392; "Q"=118.24386808843; "R"=1064.19481279587; "S"=0.1015796698791689; "T"=0.1015796698791689 }
    3 = @{ "E"=3; "G"=31.29437; "H"=93.88310999999999; "I"=0.2388439922596655; "J"=0.2388439922596655; "K"=3; "M"=4.333403333333333; "N"=13.00021; "O"=0.4877633593505858; "P"=0.4877633593505858; "Q"=135.6111272725666; "R"=1220.5001454531; "S"=0.1164993480252798; "T"=0.1164993480252798 }
    4 = @{ "E"=3; "G"=31.29437; "H"=93.88310999999999; "I"=0.2388439922596655; "J"=0.2388439922596655; "K"=3; "M"=0.2909853333333334; "N"=0.8729560000000001; "O"=0.03275300561492853; "P"=0.03275300561492853; "Q"=9.106202685906666; "R"=81.95582417316; "S"=0.007822858619572774; "T"=0.007822858619572774 }
    5 = @{ "E"=3; "G"=31.29437; "H"=93.88310999999999; "I"=0.2388439922596655; "J"=0.2388439922596655; "K"=3; "M"=0.4814053333333333; "N"=1.444216; "O"=0.0541864822020464; "P"=0.0541864822020464; "Q"=15.06527662130666; "R"=135.58748959176; "S"=0.01294211573564408; "T"=0.01294211573564408 }
    6 = @{ "E"=3; "G"=33.86972933333333; "H"=101.609188; "I"=0.2584995758255442; "J"=0.2584995758255442; "K"=3; "M"=3.778439; "N"=11.335317; "O"=0.4252971528324392; "P"=0.4252971528324392; "Q"=127.9747062325107; "R"=1151.772356092596; "S"=0.1099391336069972; "T"=0.1099391336069972 }
    7 = @{ "E"=3; "G"=33.86972933333333; "H"=101.609188; "I"=0.2584995758255442; "J"=0.2584995758255442; "K"=3; "M"=4.333403333333333; "N"=13.00021; "O"=0.4877633593505858; "P"=0.4877633593505858; "Q"=146.7711979921644; "R"=1320.94078192948; "S"=0.1260866214953689; "T"=0.1260866214953689 }
    8 = @{ "E"=3; "G"=33.86972933333333; "H"=101.609188; "I"=0.2584995758255442; "J"=0.2584995758255442; "K"=3; "M"=0.2909853333333334; "N"=0.8729560000000001; "O"=0.03275300561492853; "P"=0.03275300561492853; "Q"=9.855594479969779; "R"=88.700350319728; "S"=0.008466638058470692; "T"=0.008466638058470692 }
    9 = @{ "E"=3; "G"=33.86972933333333; "H"=101.609188; "I"=0.2584995758255442; "J"=0.2584995758255442; "K"=3; "M"=0.4814053333333333; "N"=1.444216; "O"=0.0541864822020464; "P"=0.05418648220204641; "Q"=16.30506833962311; "R"=146.745615056608; "S"=0.01400718266470739; "T"=0.01400718266470739 }
    10 = @{ "E"=3; "G"=10.670404; "H"=32.011212; "I"=0.08143835106389757; "J"=0.08143835106389757; "K"=3; "M"=3.778439; "N"=11.335317; "O"=0.4252971528324392; "P"=0.4252971528324392; "Q"=40.317470619356; "R"=362.857235574204; "S"=0.03463549883884429; "T"=0.03463549883884429 }
    11 = @{ "E"=3; "G"=10.670404; "H"=32.011212; "I"=0.08143835106389757; "J"=0.08143835106389757; "K"=3; "M"=4.333403333333333; "N"=13.00021; "O"=0.4877633593505858; "P"=0.4877633593505858; "Q"=46.23916426161333; "R"=416.15247835452; "S"=0.03972264369489904; "T"=0.03972264369489904 }
    12 = @{ "E"=3; "G"=10.670404; "H"=32.011212; "I"=0.08143835106389757; "J"=0.08143835106389757; "K"=3; "M"=0.2909853333333334; "N"=0.8729560000000001; "O"=0.03275300561492853; "P"=0.03275300561492853; "Q"=3.104931064741334; "R"=27.944379582672; "S"=0.002667350769666358; "T"=0.002667350769666358 }
    13 = @{ "E"=3; "G"=10.670404; "H"=32.011212; "I"=0.08143835106389757; "J"=0.08143835106389757; "K"=3; "M"=0.4814053333333333; "N"=1.444216; "O"=0.0541864822020464; "P"=0.05418648220204641; "Q"=5.136789394421332; "R"=46.231104549792; "S"=0.004412857760487892; "T"=0.004412857760487893 }
    14 = @{ "E"=3; "G"=55.18980966666666; "H"=165.569429; "I"=0.4212180808508926; "J"=0.4212180808508926; "K"=3; "M"=3.778439; "N"=11.335317; "O"=0.4252971528324392; "P"=0.4252971528324392; "Q"=208.5313292471103; "R"=1876.781963223993; "S"=0.1791428505074288; "T"=0.1791428505074288 }
    15 = @{ "E"=3; "G"=55.18980966666666; "H"=165.569429; "I"=0.4212180808508926; "J"=0.4212180808508926; "K"=3; "M"=4.333403333333333; "N"=13.00021; "O"=0.4877633593505858; "P"=0.4877633593505858; "Q"=239.1597051755655; "R"=2152.43734658009; "S"=0.205454746135038; "T"=0.2054547461350381 }
    16 = @{ "E"=3; "G"=55.18980966666666; "H"=165.569429; "I"=0.4212180808508926; "J"=0.4212180808508926; "K"=3; "M"=0.2909853333333334; "N"=0.8729560000000001; "O"=0.03275300561492853; "P"=0.03275300561492853; "Q"=16.05942516245822; "R"=144.534826462124; "S"=0.01379615816721871; "T"=0.01379615816721871 }
    17 = @{ "E"=3; "G"=55.18980966666666; "H"=165.569429; "I"=0.4212180808508926; "J"=0.4212180808508926; "K"=3; "M"=0.4814053333333333; "N"=1.444216; "O"=0.0541864822020464; "P"=0.05418648220204641; "Q"=26.56866871918488; "R"=239.118018472664; "S"=0.02282432604120704; "T"=0.02282432604120704 }
}

$colIndex = @{
    "E" = 5
    "G" = 7
    "H" = 8
    "I" = 9
    "J" = 10
    "K" = 11
    "M" = 13
    "N" = 14
    "O" = 15
    "P" = 16
    "Q" = 17
    "R" = 18
    "S" = 19
    "T" = 20
}

foreach ($rowNum in $rowUpdates.Keys) {
    $cols = $rowUpdates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $colNum = $colIndex[$colLetter]
        $ws.Cells.Item([int]$rowNum, $colNum).Value = $cols[$colLetter]
    }
}

Write-Output "Updated $($rowUpdates.Keys.Count) rows"